# Scheduled-runner refresh of cached market/profit figures in Sheets
# (Fenrir_Profits.xlsx). Only the numeric result columns (H..N) on a
# handful of rows change; item/recipe metadata columns A..G are untouched.

$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 5796.6665
$ws.Range("I76").Value = 4753.6
$ws.Range("J76").Value = 7100.5
$ws.Range("K76").Value = 4753.6
$ws.Range("L76").Value = 7100.5
$ws.Range("M76").Value = -4438.6
$ws.Range("N76").Value = -7730.5

$ws.Range("H79").Value = 5796.6665
$ws.Range("I79").Value = 4753.6
$ws.Range("J79").Value = 7100.5
$ws.Range("K79").Value = 4753.6
$ws.Range("L79").Value = 7100.5
$ws.Range("M79").Value = -3661.6
$ws.Range("N79").Value = -9284.5

$ws.Range("H94").Value = 9765.909
$ws.Range("I94").Value = 10672.5
$ws.Range("K94").Value = 10672.5
$ws.Range("M94").Value = -10221.5

$ws.Range("H98").Value = 2851827.8
$ws.Range("I98").Value = 3586486.5
$ws.Range("J98").Value = 5025.125
$ws.Range("K98").Value = 3586486.5
$ws.Range("L98").Value = 5025.125
$ws.Range("M98").Value = -3584988.5
$ws.Range("N98").Value = -8021.125

$ws.Range("H122").Value = 2851827.8
$ws.Range("I122").Value = 3586486.5
$ws.Range("J122").Value = 5025.125
$ws.Range("K122").Value = 10759459.5
$ws.Range("L122").Value = 15075.375
$ws.Range("M122").Value = -10757009.5
$ws.Range("N122").Value = -19975.375

$ws.Range("H135").Value = 3345.4412
$ws.Range("I135").Value = 1854.6875
$ws.Range("J135").Value = 4670.5557
$ws.Range("K135").Value = 16692.1875
$ws.Range("L135").Value = 42035.0013
$ws.Range("M135").Value = -14157.1875
$ws.Range("N135").Value = -47105.0013

$ws.Range("H138").Value = 1837.98
$ws.Range("I138").Value = 999.7368
$ws.Range("J138").Value = 2351.742
$ws.Range("K138").Value = 2999.2104
$ws.Range("L138").Value = 7055.226000000001
$ws.Range("M138").Value = 2140.7896
$ws.Range("N138").Value = -17335.226

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11366.57
$ws.Range("I32").Value = 9012.727000000001
$ws.Range("J32").Value = 23724.25
$ws.Range("K32").Value = 9012.727000000001
$ws.Range("L32").Value = 23724.25
$ws.Range("M32").Value = -8725.727000000001
$ws.Range("N32").Value = -24298.25

$ws.Range("H61").Value = 6505.7085
$ws.Range("I61").Value = 7577.5
$ws.Range("J61").Value = 3290.3333
$ws.Range("K61").Value = 7577.5
$ws.Range("L61").Value = 3290.3333
$ws.Range("M61").Value = -7365.5
$ws.Range("N61").Value = -3714.3333

$ws.Range("H74").Value = 707.10345
$ws.Range("I74").Value = 615.9048
$ws.Range("J74").Value = 946.5
$ws.Range("K74").Value = 615.9048
$ws.Range("L74").Value = 946.5
$ws.Range("M74").Value = 258.0952
$ws.Range("N74").Value = -2694.5

$ws.Range("H77").Value = 707.10345
$ws.Range("I77").Value = 615.9048
$ws.Range("J77").Value = 946.5
$ws.Range("K77").Value = 3079.524
$ws.Range("L77").Value = 4732.5
$ws.Range("M77").Value = 1288.476
$ws.Range("N77").Value = -13468.5

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 14516239
$ws.Range("I134").Value = 15175613
$ws.Range("J134").Value = 10000
$ws.Range("K134").Value = 45526839
$ws.Range("L134").Value = 30000
$ws.Range("M134").Value = -45524304
$ws.Range("N134").Value = -35070

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5669.9854
$ws.Range("I31").Value = 1150.1177
$ws.Range("J31").Value = 10060.714
$ws.Range("K31").Value = 1150.1177
$ws.Range("L31").Value = 10060.714
$ws.Range("M31").Value = -855.1177
$ws.Range("N31").Value = -10650.714

$ws.Range("H34").Value = 5669.9854
$ws.Range("I34").Value = 1150.1177
$ws.Range("J34").Value = 10060.714
$ws.Range("K34").Value = 1150.1177
$ws.Range("L34").Value = 10060.714
$ws.Range("M34").Value = -948.1177
$ws.Range("N34").Value = -10464.714

$ws.Range("H58").Value = 3599990.5
$ws.Range("I58").Value = 5329055
$ws.Range("J58").Value = 8856
$ws.Range("K58").Value = 5329055
$ws.Range("L58").Value = 8856
$ws.Range("M58").Value = -5328852
$ws.Range("N58").Value = -9262

$ws.Range("H132").Value = 6947973
$ws.Range("I132").Value = 13889820
$ws.Range("J132").Value = 6125.75
$ws.Range("K132").Value = 41669460
$ws.Range("L132").Value = 18377.25
$ws.Range("M132").Value = -41666930
$ws.Range("N132").Value = -23437.25

$ws.Range("H134").Value = 14206618
$ws.Range("I134").Value = 25002772
$ws.Range("J134").Value = 5209823
$ws.Range("K134").Value = 75008316
$ws.Range("L134").Value = 15629469
$ws.Range("M134").Value = -75005781
$ws.Range("N134").Value = -15634539

$ws.Range("H136").Value = 3599990.5
$ws.Range("I136").Value = 5329055
$ws.Range("J136").Value = 8856
$ws.Range("K136").Value = 15987165
$ws.Range("L136").Value = 26568
$ws.Range("M136").Value = -15984615
$ws.Range("N136").Value = -31668

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 639.1957
$ws.Range("J122").Value = 954.2
$ws.Range("L122").Value = 8587.800000000001
$ws.Range("N122").Value = -13487.8

# --- GSM ---
# Row 44: H/J/L reset to 0 (I/K already 0, unchanged) and the N cell is
# dropped entirely (was -6192, no replacement value now -- matches the
# sparse M/N pattern used elsewhere in this sheet for "no profit data").
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 76924080
$ws.Range("I126").Value = 987.9
$ws.Range("J126").Value = 333334400
$ws.Range("K126").Value = 2963.7
$ws.Range("L126").Value = 1000003200
$ws.Range("M126").Value = -493.6999999999998
$ws.Range("N126").Value = -1000008140

$ws.Range("H132").Value = 682667.6
$ws.Range("I132").Value = 58115.223
$ws.Range("J132").Value = 3493153.5
$ws.Range("K132").Value = 174345.669
$ws.Range("L132").Value = 10479460.5
$ws.Range("M132").Value = -171815.669
$ws.Range("N132").Value = -10484520.5
